# Goterra Budget workbook update
# - Fix a typo in F3
# - Highlight the construction-cost labels (A6:A13) with a fill colour
# - Turn several reference URLs into real hyperlinks with matching style
# - Insert a new row (24) so the "Pay Guide" / "Section" notes sit on their
#   own rows, and re-word them as bold-prefixed rich text
# - Add a new footnote in F4 and a new sourced note in row 26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix "Variabe" -> "Variable" typo in F3
$ws.Range("F3").Value = "*** Variable depending on Goterras Needs"

# 2. Add the accent background colour used elsewhere in the sheet (same as
#    B2/E2) to the construction-cost row labels A6:A13, keeping their
#    existing borders.
$ws.Range("A6:A13").Interior.Color = $ws.Range("B2").Interior.Color

# 3. Insert a new row at 24 - this pushes the old row24 ("5" / canstarblue
#    link) down to row25, and old row25 ("6" / Sensitar text) down to row26.
$ws.Rows(24).Insert()

# --- Row 20: Water Tank source link ---
$ws.Hyperlinks.Add($ws.Range("B20"), "https://www.alibaba.com/product-detail/FRP-SMC-GRP-Stainless-steel-Galvanized_60815777657.html?spm=a2700.7724838.2017115.25.1998399amtiN7A&s=p") | Out-Null

# --- Row 21: Cooling Tower source link ---
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.alibaba.com/product-detail/Small-Cooling-Tower-Industrial-Industry-Water_875962686.html?spm=a2700.7724838.2017115.12.2ac846425fFSmm") | Out-Null

# --- Row 23: Pay Guide note -> bold "Pay Guide - Meat Industry Award 2010:" + link text ---
$ws.Range("E23").ClearContents()
$ws.Range("B23").Value = "Pay Guide - Meat Industry Award 2010: https://github.com/JessYJY/InsectFarming/blob/master/Research/meat-industry-award-ma000059-pay-guide.pdf"
$ws.Range("B23").Characters(1, 37).Font.Bold = $true

# --- Row 24 (new, blank row): Section note -> bold "Section:" + description ---
$ws.Range("B24").Value = "Section: Adult-Meat processing establishment - Full-time & part-time -level 5"
$ws.Range("B24").Characters(1, 8).Font.Bold = $true

# --- Row 25 (was row24): now the canstarblue electricity source link ---
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "https://www.canstarblue.com.au/electricity/canberra-act-electricity/"
$ws.Hyperlinks.Add($ws.Range("B25"), "https://www.canstarblue.com.au/electricity/canberra-act-electricity/") | Out-Null

# --- Row 26 (was row25): now the Sensitar equipment note, with row-level fill ---
$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Sensitar Equipment, PDF of equipment and prices listed in repositry: https://github.com/JessYJY/InsectFarming/blob/master/Communication/Expert/Sensitar/Sensitar%201T%20rendering%20plant%20price.pdf "
$ws.Rows(26).Interior.Color = $ws.Range("B2").Interior.Color

# 4. New footnote next to the power-consumption calc
$ws.Range("F4").Value = "*** Goterra would be able to provide a better estimate for this, this is currently based on residentila postcodes"

# 5. Update the active selection to match the author's final cursor position
$ws.Range("C13").Select()
